# Nalco aluminium-ingot price sheet: refresh the 30 data rows (rows 2-31)
# that referenced the 01-11-2025 circular with the newer 02-11-2025 circular:
#   - Basic Price (D): 297.15 -> 296.05
#   - Circular Date (E): 01-11-2025 -> 02-11-2025
#   - Circular Link (F): ...Ingot-01-11-2025.pdf -> ...Ingot-02-11-2025.pdf
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 4).Value = 296.05

    # Force the text number format first so the dd-mm-yyyy-looking string is
    # stored as literal text instead of being auto-converted to a date serial.
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = "02-11-2025"

    $ws.Cells.Item($r, 6).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"
}
